$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Weapons")

# Capture current (pre-edit) row 28 ("Nordic") and row 29 ("Orcish") values across columns A:H
$row28 = @{}
$row29 = @{}
for ($c = 1; $c -le 8; $c++) {
    $row28[$c] = $ws.Cells.Item(28, $c).Value2
    $row29[$c] = $ws.Cells.Item(29, $c).Value2
}

# The "Nordic" row (28) is being renamed to "Quicksilver"; since the weapon list is kept
# alphabetically sorted, it now sorts after "Orcish" (row 29) instead of before it, so the
# two rows swap places. First clear both rows, then write back the new contents.
for ($c = 1; $c -le 8; $c++) {
    $ws.Cells.Item(28, $c).ClearContents()
    $ws.Cells.Item(29, $c).ClearContents()
}

# New row 28 = old "Orcish" row (unchanged content), moved up from row 29
for ($c = 1; $c -le 8; $c++) {
    $v = $row29[$c]
    if ($null -ne $v) {
        $ws.Cells.Item(28, $c).Value = $v
    }
}

# New row 29 = old "Nordic" row, moved down from row 28, with its name changed to "Quicksilver"
for ($c = 1; $c -le 8; $c++) {
    $v = $row28[$c]
    if ($null -ne $v) {
        $ws.Cells.Item(29, $c).Value = $v
    }
}
$ws.Cells.Item(29, 1).Value = "Quicksilver"

# --- Selection / active-tab bookkeeping to match the authored workbook view state ---

# Weapons sheet becomes the active tab, with A28 selected
$ws.Activate()
$ws.Range("A28").Select()

# Artifacts sheet is no longer the tab that was last active before saving
$wsArtifacts = $wb.Worksheets.Item("Artifacts")
$wsArtifacts.Range("B58").Select()
$ws.Activate()
